# Apply the upstream "Automatic update of files" sync: the source export
# re-ordered/re-fetched observation records, so rows 3-6 need their field
# values updated in place to match the new export snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (was "Revlummer" 112043158) -> "Tretåig hackspett" record, new Id,
# and the public comment is cleared in the new snapshot.
$ws.Range("A3").Value = 112042452
$ws.Range("B3").Value = 56398
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("P3").Value = "Stor Mpmerg, Kilen-Stor, Moberg, Leksand, Dlr"
$ws.Range("Q3").Value = 511613.7990622812
$ws.Range("R3").Value = 6733639.811082688
$ws.Range("S3").Value = 25
$ws.Range("Z3").Value = "00:00"
$ws.Range("AB3").Value = "00:00"
$ws.Range("AC3").Value = ""
$ws.Range("AW3").Value = "Åke Sköld"
$ws.Range("AX3").Value = "Åke Sköld"

# Row 4 (was "Tretåig hackspett" 112044333) -> "Blåsippa" record.
$ws.Range("A4").Value = 112042940
$ws.Range("B4").Value = 98535
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 222498
$ws.Range("F4").Value = "Blåsippa"
$ws.Range("G4").Value = "Hepatica nobilis"
$ws.Range("H4").Value = "Schreb."
$ws.Range("P4").Value = "Stor-Moberg (Stor-Moberg), Dlr"
$ws.Range("Q4").Value = 511610.9043343531
$ws.Range("R4").Value = 6733626.107665217
$ws.Range("S4").Value = 1
$ws.Range("Z4").Value = "10:33"
$ws.Range("AB4").Value = "10:33"
$ws.Range("AC4").Value = "Fullt med blåsippsblad på denna sidan bäcken"
$ws.Range("AW4").Value = "Evalena Sköld"
$ws.Range("AX4").Value = "Evalena Sköld, Åke Sköld"

# Row 5 (was "Svavelriska" 112043031) -> "Revlummer" record.
$ws.Range("A5").Value = 112043158
$ws.Range("B5").Value = 95532
$ws.Range("E5").Value = 221945
$ws.Range("F5").Value = "Revlummer"
$ws.Range("G5").Value = "Lycopodium annotinum"
$ws.Range("H5").Value = "L."
$ws.Range("Q5").Value = 511628.0588172724
$ws.Range("R5").Value = 6733623.228879539
$ws.Range("Z5").Value = "10:51"
$ws.Range("AB5").Value = "10:51"
$ws.Range("AC5").Value = "Finns fläckvis i området"
$ws.Range("AX5").Value = "Evalena Sköld, Åke Sköld"

# Row 6 (was "Blåsippa" 112042940) -> "Svavelriska" record.
$ws.Range("A6").Value = 112043031
$ws.Range("B6").Value = 90332
$ws.Range("E6").Value = 4769
$ws.Range("F6").Value = "Svavelriska"
$ws.Range("G6").Value = "Lactarius scrobiculatus"
$ws.Range("H6").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q6").Value = 511625.1419049087
$ws.Range("R6").Value = 6733616.372369035
$ws.Range("Z6").Value = "10:42"
$ws.Range("AB6").Value = "10:42"
$ws.Range("AC6").Value = ""
$ws.Range("AX6").Value = "Evalena Sköld"
